$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 218-220 (column B and D values revised) ---
$ws.Cells.Item(218, 2).Value = 603042000000
$ws.Cells.Item(218, 4).Value = 131734713939.3145

$ws.Cells.Item(219, 2).Value = 603647900000
$ws.Cells.Item(219, 4).Value = 133394007027.1584

$ws.Cells.Item(220, 2).Value = 612183900000
$ws.Cells.Item(220, 4).Value = 131536473217.1634

# --- Append new row 224 ---
# Copy formatting (date style) from the last existing data row (A223) to A224
$ws.Cells.Item(223, 1).Copy($ws.Cells.Item(224, 1))

$ws.Cells.Item(224, 1).Value = 45078
$ws.Cells.Item(224, 2).Value = 624519300000
$ws.Cells.Item(224, 3).Value = 0.2204342554833021
$ws.Cells.Item(224, 4).Value = 137665446930.453

Write-Host "Edit applied successfully"
